$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(61.22999954223633, 64.70999908447266, 70.48999786376953, 59.5099983215332, 55696909)
    3 = @(72.16999816894531, 74.11000061035156, 75.88999938964844, 70.51000213623047, 55696909)
    4 = @(74.23999786376953, 77.34999847412109, 81.94999694824219, 73.22000122070312, 55696909)
    5 = @(77.97000122070312, 74.90000152587891, 78.04000091552734, 64.81999969482422, 55696909)
    6 = @(74.19000244140625, 72.93000030517578, 78.40000152587891, 72.11000061035156, 55696909)
    7 = @(64.38999938964844, 70.23999786376953, 71.08999633789062, 61.47000122070312, 55696909)
    8 = @(69.19999694824219, 64.37000274658203, 69.19999694824219, 62.95000076293945, 55696909)
    9 = @(64.68000030517578, 64.36000061035156, 67.56999969482422, 63.31999969482422, 55696909)
    10 = @(75.31999969482422, 77, 78.48999786376953, 73.48999786376953, 55696909)
    11 = @(84.69999694824219, 85.93000030517578, 87.93000030517578, 81.56999969482422, 55696909)
    12 = @(87.94000244140625, 91.15000152587891, 91.97000122070312, 86.52999877929688, 55696909)
    13 = @(107.620002746582, 117.4800033569336, 119.1800003051758, 106.9000015258789, 55696909)
    14 = @(114.3499984741211, 114.3499984741211, 119.7200012207031, 110.1999969482422, 55696909)
    15 = @(123.7099990844727, 130.2100067138672, 137.7299957275391, 121.8199996948242, 55696909)
    16 = @(138.3099975585938, 119.4700012207031, 139.3099975585938, 110.7200012207031, 55696909)
    17 = @(113.7900009155273, 141.4799957275391, 142.2599945068359, 111.4400024414062, 55696909)
    18 = @(170.4700012207031, 179.3600006103516, 179.9299926757812, 166.5700073242188, 55696909)
    19 = @(175.9100036621094, 193.7899932861328, 201, 174.8800048828125, 55696909)
    20 = @(183.0700073242188, 175.9600067138672, 191.6399993896484, 168.2599945068359, 55696909)
    21 = @(214.1199951171875, 228.1399993896484, 236.0399932861328, 212, 55696909)
    22 = @(179.1000061035156, 220.8899993896484, 226.0500030517578, 158.8300018310547, 55696909)
    23 = @(252.75, 290.0799865722656, 290.1199951171875, 246.009994506836, 55696909)
    24 = @(328.3200073242188, 308.9500122070312, 356.4800109863281, 306.2999877929688, 55696909)
    25 = @(359.010009765625, 344.4299926757812, 369.510009765625, 328.8999938964844, 55696909)
    26 = @(401.6300048828125, 457.75, 461.9800109863281, 394.0199890136719, 55696909)
    27 = @(510.4299926757813, 559.7999877929688, 563.8900146484375, 507.7999877929688, 55696909)
    28 = @(574.0800170898438, 673.239990234375, 674.8099975585938, 552.0800170898438, 55696909)
    29 = @(670.3200073242188, 476.1400146484375, 675.3400268554688, 423.5400085449219, 55696909)
    30 = @(297.3699951171875, 264.989990234375, 321.2099914550781, 260.6799926757812, 55696909)
    31 = @(296.6499938964844, 349.25, 349.6700134277344, 291.4299926757812, 55696909)
    32 = @(365.8599853515625, 350, 389.5899963378906, 303.5499877929688, 55696909)
    33 = @(334.7799987792969, 332.6499938964844, 371.1000061035156, 301.3699951171875, 55696909)
    34 = @(294.9299926757812, 282.4400024414062, 304.3299865722656, 271.0400085449219, 55696909)
    35 = @(224.979995727539, 236.8099975585937, 255.1900024414062, 224.5500030517578, 55696909)
    36 = @(255.2200012207031, 217.5700073242188, 256.25, 209.1499938964844, 55696909)
    37 = @(293.489990234375, 278.1099853515625, 309.3800048828125, 276.2000122070312, 55696909)
    38 = @(275.5599975585937, 235.259994506836, 277.0400085449219, 234.8200073242188, 55696909)
    39 = @(186.4199981689453, 215.1300048828125, 223.0299987792969, 179.5500030517578, 55696909)
    40 = @(199.1900024414062, 188.6499938964844, 204.5500030517578, 188.1000061035156, 55696909)
    41 = @(236.1399993896484, 253.9600067138672, 257.2300109863281, 222.5, 55696909)
    42 = @(169.1999969482422, 156.9100036621094, 172.2100067138672, 138.1499938964844, 55696909)
    43 = @(176.6300048828125, 157.7100067138672, 185.8200073242188, 157.2700042724609, 55696909)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $vals[1]
    $ws.Cells.Item($r, 6).Value = $vals[2]
    $ws.Cells.Item($r, 7).Value = $vals[3]
    $ws.Cells.Item($r, 8).Value = $vals[4]
    $ws.Cells.Item($r, 9).Value = "EPAM"
}
